$d = $word.ActiveDocument
$d.Content.Find.Execute("Asddinfg", $true, $false, $false, $false, $false, $true, 1, $false, "Asddinfg Making nother change", 2)
